# Inserts the new "entities and entity descriptions" paragraphs for the
# Lottery data model Readme, right after the second (empty, centered)
# paragraph and before the trailing empty paragraph / sectPr.
$d = $word.ActiveDocument

$newTexts = @(
    'Below are the entities and entity descriptions for the Lottery data model. ',
    'Games:',
    'Games have the following attributes; GameId (primary key, identity), GameName, Rules, GameDescription, and Cost. ',
    'States have the following attributes; StateId (primary key, identity), StateName',
    'Games and States are related with the associative table, StateToGame. StateToGame has a unique identifier, Id (primary key, identity), as well as a column for StateId and GameId, which reference the appropriate columns from the Games and States entities.',
    'Games have odds associated to them, which users can read as string values. The Odds entity has the following attributes; GameOddsId (primary key, identity), GameID (foreign key of Games), Match, Prize, and Odds (all strings describing the amount of numbers matched, the prize, and the odds).',
    'All Games have drawings which occur on various dates. The Drawing entity stores date information, jackpot value, and game name with the following attirbutes; DrawingId (primary key, id), DrawingDate, Jackpot, and GameId (foreign key of Games).',
    'All Drawings consist of various numbers. The Number entity represents a single number in a drawing, associated with the drawing by referencing the specific DrawingId. It has the following attributes; NumberId (primary key, identity), Number, NumberTypeId (foreign key of NumberType, described below), and DrawingId (foreign key of Drawing).',
    'Lastly, the NumberType table is a lookup table for the type of number in the Number table. For example a normal ball (single number in a drawing) would have the type “normal” under NumberTypeId of 1. The attributed are as follows; NumberTypeId (primary key, identity) and NumberType (the string describing the number type). ',
    'The SQL code included in this repo has all the information needed to generate this database and all its relationships. '
)

# 0-based index (within $newTexts) of the paragraph that gets a bottom border
# ("Lastly, the NumberType table is a lookup table ...").
$borderParaOffset = 8

# Anchor: the second paragraph in the original document (empty, centered).
$anchorIndex = 2

# Pass 1: create every paragraph and set its text + alignment. Paragraph
# (pPr-level) formatting such as borders must NOT be touched yet, because
# Word carries the "current" paragraph's mark formatting forward onto the
# *next* InsertParagraphAfter() call — setting the border before the loop
# is done would leak it onto every paragraph inserted afterwards.
$borderParaIndex = -1
for ($i = 0; $i -lt $newTexts.Count; $i++) {
    $anchorRange = $d.Paragraphs($anchorIndex).Range
    $anchorRange.InsertParagraphAfter()
    $anchorIndex = $anchorIndex + 1
    $newPara = $d.Paragraphs($anchorIndex)
    $newPara.Range.Text = $newTexts[$i]
    # New paragraphs are plain (left-aligned, no pPr/jc) unlike the
    # centered title paragraphs above them.
    $newPara.Alignment = 0

    if ($i -eq $borderParaOffset) {
        $borderParaIndex = $anchorIndex
    }
}

# Pass 2: now that no more paragraphs will be inserted, it is safe to add
# the bottom border to the "Lastly, the NumberType ..." paragraph only.
$borderPara = $d.Paragraphs($borderParaIndex)
$borderPara.Borders.DistanceFromBottom = 1
$bottomBorder = $borderPara.Borders(-3)
$bottomBorder.LineStyle = 1
$bottomBorder.LineWidth = 3
$bottomBorder.Color = -16777216

Write-Output "Inserted $($newTexts.Count) paragraphs; document now has $($d.Paragraphs.Count) paragraphs."
